# Update the UUID-based file names and timestamps that were regenerated
# when the localization status report was produced again for handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "6309d183-d698-485f-9487-52c536b19194"
$newGuid = "1971f2fe-84d0-4b18-b839-73772f7545b3"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2017-02-09 11:02:41"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.7ed5c7830191e80c554ec1ee933854866cac2019.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-02-09 11:02:18"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.7ed5c7830191e80c554ec1ee933854866cac2019.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-02-09 11:02:41"

# Keep the hyperlink display text for B2/A2 in sync with the new file names.
# NOTE: must mutate the Hyperlink objects returned from enumeration (foreach)
# rather than via the Hyperlinks.Item(...) indexer, since re-assigning
# TextToDisplay through Item() creates a duplicate hyperlink instead of
# updating the existing one in place.
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
